$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final consensus clustering definition: the moffitt_stromal_type column (G)
# no longer carries a stroma call ("Normal Stroma" / "Activated Stroma" / "NA")
# for any sample - every data row's value is cleared to an empty string.
#
# Assign an empty text value (leading "'" forces a genuine empty-text cell
# rather than dropping the cell to a blank/number cell) to every row.
for ($r = 2; $r -le 37; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
